$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-18 07:48:59'
$ws.Range('E3').Value = '2026-02-18 07:49:02'
$ws.Range('K3').Value = '0.0 MJ/m2'
$ws.Range('M3').Value = '-2.0 °C 7:27 TU'
$ws.Range('O3').Value = '-3.5 °C'
$ws.Range('E4').Value = '2026-02-18 07:49:05'
$ws.Range('J4').Value = '1018.1 hPa'
$ws.Range('K4').Value = '0.0 MJ/m2'
$ws.Range('N4').Value = '4.2 °C 7:02 TU'
$ws.Range('O4').Value = '6.6 °C'
$ws.Range('E5').Value = '2026-02-18 07:49:07'
$ws.Range('H5').NumberFormat = '@'
$ws.Range('H5').Value = '83%'
$ws.Range('K5').Value = '0.0 MJ/m2'
$ws.Range('M5').Value = '0.9 °C 7:23 TU'
$ws.Range('O5').Value = '-1.5 °C'
$ws.Range('E6').Value = '2026-02-18 07:49:10'
$ws.Range('K6').Value = '0.1 MJ/m2'
$ws.Range('N6').Value = '5.9 °C 7:04 TU'
$ws.Range('O6').Value = '7.6 °C'
$ws.Range('E7').Value = '2026-02-18 07:49:12'
$ws.Range('J7').Value = '1018.0 hPa'
$ws.Range('K7').Value = '0.1 MJ/m2'
$ws.Range('E8').Value = '2026-02-18 07:49:15'
$ws.Range('J8').Value = '1018.1 hPa'
$ws.Range('K8').Value = '0.1 MJ/m2'
$ws.Range('L8').Value = '19.1 km/h - 254º 7:19 TU'
$ws.Range('N8').Value = '7.9 °C 7:18 TU'
$ws.Range('O8').Value = '8.6 °C'
$ws.Range('E9').Value = '2026-02-18 07:49:17'
$ws.Range('K9').Value = '0.1 MJ/m2'
$ws.Range('E10').Value = '2026-02-18 07:49:20'
$ws.Range('K10').Value = '0.1 MJ/m2'
$ws.Range('O10').Value = '6.3 °C'
$ws.Range('E11').Value = '2026-02-18 07:49:22'
$ws.Range('N11').Value = '-0.4 °C 7:06 TU'
$ws.Range('O11').Value = '1.5 °C'
$ws.Range('E12').Value = '2026-02-18 07:49:25'
$ws.Range('E13').Value = '2026-02-18 07:49:27'
$ws.Range('J13').Value = '1023.0 hPa'
$ws.Range('K13').Value = '0.0 MJ/m2'
$ws.Range('O13').Value = '-2.6 °C'
$ws.Range('E14').Value = '2026-02-18 07:49:30'
$ws.Range('K14').Value = '0.1 MJ/m2'
$ws.Range('L14').Value = '12.6 km/h - 309º 7:27 TU'
$ws.Range('O14').Value = '9.7 °C'
$ws.Range('E15').Value = '2026-02-18 07:49:32'
$ws.Range('E16').Value = '2026-02-18 07:49:35'
$ws.Range('H16').NumberFormat = '@'
$ws.Range('H16').Value = '37%'
$ws.Range('K16').Value = '0.0 MJ/m2'
$ws.Range('O16').Value = '0.7 °C'
$ws.Range('E17').Value = '2026-02-18 07:49:37'
$ws.Range('K17').Value = '0.1 MJ/m2'
$ws.Range('E18').Value = '2026-02-18 07:49:40'
$ws.Range('J18').Value = '1018.1 hPa'
$ws.Range('K18').Value = '0.1 MJ/m2'
$ws.Range('O18').Value = '7.0 °C'
$ws.Range('E19').Value = '2026-02-18 07:49:43'
$ws.Range('N19').Value = '5.1 °C 7:12 TU'
$ws.Range('E20').Value = '2026-02-18 07:49:45'
$ws.Range('H20').NumberFormat = '@'
$ws.Range('H20').Value = '74%'
$ws.Range('K20').Value = '0.0 MJ/m2'
$ws.Range('M20').Value = '0.6 °C 7:28 TU'
$ws.Range('O20').Value = '-0.9 °C'
$ws.Range('E21').Value = '2026-02-18 07:49:47'
$ws.Range('H21').NumberFormat = '@'
$ws.Range('H21').Value = '86%'
$ws.Range('J21').Value = '1020.6 hPa'
$ws.Range('K21').Value = '0.0 MJ/m2'
$ws.Range('N21').Value = '0.1 °C 7:00 TU'
$ws.Range('E22').Value = '2026-02-18 07:49:50'
$ws.Range('E23').Value = '2026-02-18 07:49:52'
$ws.Range('K23').Value = '0.0 MJ/m2'
$ws.Range('O23').Value = '0.7 °C'
$ws.Range('E24').Value = '2026-02-18 07:49:55'
$ws.Range('J24').Value = '1018.7 hPa'
$ws.Range('O24').Value = '4.7 °C'
$ws.Range('E25').Value = '2026-02-18 07:49:57'
$ws.Range('H25').NumberFormat = '@'
$ws.Range('H25').Value = '42%'
$ws.Range('K25').Value = '0.0 MJ/m2'
$ws.Range('M25').Value = '2.4 °C 7:29 TU'
$ws.Range('O25').Value = '0.0 °C'
$ws.Range('E26').Value = '2026-02-18 07:50:00'
$ws.Range('E27').Value = '2026-02-18 07:50:02'
$ws.Range('H27').NumberFormat = '@'
$ws.Range('H27').Value = '48%'
$ws.Range('O27').Value = '1.0 °C'
$ws.Range('E28').Value = '2026-02-18 07:50:05'
$ws.Range('O28').Value = '4.4 °C'
$ws.Range('E29').Value = '2026-02-18 07:50:07'
$ws.Range('H29').NumberFormat = '@'
$ws.Range('H29').Value = '93%'
$ws.Range('E30').Value = '2026-02-18 07:50:10'
$ws.Range('J30').Value = '1018.1 hPa'
$ws.Range('K30').Value = '0.1 MJ/m2'
$ws.Range('E31').Value = '2026-02-18 07:50:12'
$ws.Range('J31').Value = '1016.7 hPa'
$ws.Range('K31').Value = '0.0 MJ/m2'
$ws.Range('N31').Value = '9.4 °C 7:06 TU'
$ws.Range('E32').Value = '2026-02-18 07:50:15'
$ws.Range('I32').Value = '0.1 mm'
$ws.Range('K32').Value = '0.0 MJ/m2'
$ws.Range('E33').Value = '2026-02-18 07:50:18'
$ws.Range('O33').Value = '-0.3 °C'
$ws.Range('E34').Value = '2026-02-18 07:50:20'
$ws.Range('O34').Value = '0.7 °C'
$ws.Range('E35').Value = '2026-02-18 07:50:23'
$ws.Range('J35').Value = '1019.0 hPa'
$ws.Range('K35').Value = '0.0 MJ/m2'
$ws.Range('O35').Value = '6.7 °C'
$ws.Range('E36').Value = '2026-02-18 07:50:25'
$ws.Range('K36').Value = '0.1 MJ/m2'
$ws.Range('E37').Value = '2026-02-18 07:50:28'
$ws.Range('J37').Value = '1021.0 hPa'
$ws.Range('O37').Value = '1.0 °C'
$ws.Range('E38').Value = '2026-02-18 07:50:30'
$ws.Range('H38').NumberFormat = '@'
$ws.Range('H38').Value = '95%'
$ws.Range('K38').Value = '0.1 MJ/m2'
$ws.Range('O38').Value = '8.1 °C'
$ws.Range('E39').Value = '2026-02-18 07:50:33'
$ws.Range('K39').Value = '0.0 MJ/m2'
$ws.Range('E40').Value = '2026-02-18 07:50:35'
$ws.Range('O40').Value = '0.4 °C'
$ws.Range('E41').Value = '2026-02-18 07:50:38'
$ws.Range('E42').Value = '2026-02-18 07:50:40'
$ws.Range('O42').Value = '7.7 °C'
$ws.Range('E43').Value = '2026-02-18 07:50:43'
$ws.Range('K43').Value = '0.1 MJ/m2'
$ws.Range('O43').Value = '6.6 °C'
$ws.Range('E44').Value = '2026-02-18 07:50:45'
$ws.Range('H44').NumberFormat = '@'
$ws.Range('H44').Value = '62%'
$ws.Range('O44').Value = '-3.3 °C'
$ws.Range('E45').Value = '2026-02-18 07:50:47'
$ws.Range('H45').NumberFormat = '@'
$ws.Range('H45').Value = '94%'
$ws.Range('J45').Value = '1020.0 hPa'
$ws.Range('K45').Value = '0.0 MJ/m2'
$ws.Range('L45').Value = '12.2 km/h - 129º 7:12 TU'
$ws.Range('M45').Value = '3.0 °C 7:17 TU'
$ws.Range('O45').Value = '0.9 °C'
$ws.Range('E46').Value = '2026-02-18 07:50:50'
$ws.Range('K46').Value = '0.0 MJ/m2'
$ws.Range('O46').Value = '5.9 °C'
